$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.432.25"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.804.59"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'225.23"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'0.585"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'38.07"
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("D9").Value = "'0.286"
$ws.Range("E9").Value = "  -5.16%  "
$ws.Range("D10").Value = "'0.0668"
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("D11").Value = "'0.0972"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "2.067.49"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'11.04"
$ws.Range("E13").Value = "  -6.44%  "
$ws.Range("D14").Value = "1.795.95"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "34.431.55"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'0.625"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "'4.38"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").Value = "'67.74"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "'241.29"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").Value = "0.0₃0764"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "'170.07"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'7.68"
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("D27").Value = "'17.38"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "'3.74"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").Value = "'0.0511"
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "1.326.65"
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.635"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("E39").Value = "  -7.08%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").Value = "'81.41"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "'2.78"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").Value = "'0.939"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "'13.51"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "1.967.83"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'101.55"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  -5.18%  "
